# Update the "Create Location" automation sample rows on the BasicInfo
# sheet with a fresh set of location data: a brand new reference code,
# automation id and coordinates for the first row, and a new reference
# code / automation id for the second row (its coordinates stay the
# same). Finally leave the selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)

    # These columns store numeric-looking reference codes / coordinates
    # as plain text. Assigning a bare numeric-looking string would make
    # Excel silently convert the cell to a number, so prefix it with an
    # apostrophe to force a text entry, then clear the resulting
    # quote-prefix style again so no stray cell formatting is left
    # behind.
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 2 - location gets a brand new reference code, automation id and
# coordinates.
Set-TextValue $ws.Range("V2") "9000230592"
Set-TextValue $ws.Range("W2") "Automation1597929055751"
Set-TextValue $ws.Range("X2") "43.613122"
Set-TextValue $ws.Range("Y2") "-79.556162"

# Row 3 - new reference code / automation id, coordinates unchanged.
Set-TextValue $ws.Range("V3") "9000230593"
Set-TextValue $ws.Range("W3") "Automation1597929225806"

# Leave the selection where the author left it.
$ws.Range("V10").Select() | Out-Null
